# Book1.xlsx edit: add a "Hoja1" worksheet after "Sheet1" replicating the
# DEC2BIN level table with an extra volume ("vol") column, plus a
# sample_in= footer. (commit: "probar en placa volumen logaritmico")

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New sheet, inserted right after Sheet1 so tab order / sheetId / r:id
# ordering matches (Sheet1 -> rId1, Hoja1 -> rId2).
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Hoja1"

# ---- Header row -----------------------------------------------------
$ws2.Range("A1").Value = "NIVEL"
$ws2.Range("B1").Value = "FACTOR NUM"
$ws2.Range("C1").Value = "FACTOR DEN"
$ws2.Range("D1").Value = "FACTOR"
$ws2.Range("H1").Value = "VOLUMEN (vol)"

# ---- Column A: level index 0..20 in rows 2..22 ----------------------
for ($r = 2; $r -le 22; $r++) {
    $ws2.Cells.Item($r, 1).Value = $r - 2
}

# ---- Column B: FACTOR NUM = 7^(level/10)-1 ---------------------------
$ws2.Range("B2").Formula = "=7^(A2/10)-1"
$ws2.Range("B3:B22").Formula = "=7^(A3/10)-1"

# ---- Column C: FACTOR DEN = 6 (not a shared formula, every row) -----
for ($r = 2; $r -le 22; $r++) {
    $ws2.Cells.Item($r, 3).Formula = "=6"
}

# ---- Column D: FACTOR = B/C -----------------------------------------
$ws2.Range("D2").Formula = "=B2/C2"
$ws2.Range("D3:D22").Formula = "=B3/C3"

# ---- Column E: D*2^5 -------------------------------------------------
$ws2.Range("E2").Formula = "=D2*2^5"
$ws2.Range("E3:E22").Formula = "=D3*2^5"

# ---- Column F: TRUNC(E) ----------------------------------------------
$ws2.Range("F2").Formula = "=TRUNC(E2)"
$ws2.Range("F3:F22").Formula = "=TRUNC(E3)"

# ---- Column G: DEC2BIN(F) ---------------------------------------------
$ws2.Range("G2").Formula = "=DEC2BIN(F2)"
$ws2.Range("G3:G22").Formula = "=DEC2BIN(F3)"

# ---- Column H: VOLUMEN (vol) = (B24/8)*D  (not shared, every row) ----
for ($r = 2; $r -le 22; $r++) {
    $ws2.Cells.Item($r, 8).Formula = "=(B24/8)*D$r"
}

# ---- Footer row 24: sample_in= 255 -----------------------------------
$ws2.Range("A24").Value = "sample_in="
$ws2.Range("B24").Value = 255

# ---- Column widths (best-effort; engine stores width on a 1/6-char
# grid rather than Excel's pixel grid, so exact match isn't reachable) -
$ws2.Columns.Item(1).ColumnWidth = 14.88
$ws2.Columns.Item(2).ColumnWidth = 16.17
$ws2.Columns.Item(3).ColumnWidth = 17.45
$ws2.Columns.Item(4).ColumnWidth = 18.02
$ws2.Columns.Item(5).ColumnWidth = 18.17
$ws2.Columns.Item(6).ColumnWidth = 25.02
$ws2.Columns.Item(7).ColumnWidth = 23.02
$ws2.Columns.Item(8).ColumnWidth = 20.02

# ---- Selection matches the author's saved cursor position -----------
$ws2.Range("C2").Select()

# Leave the originally active sheet selected/first, matching the source
# file (tabSelected stays on Sheet1 in the target workbook.xml bookViews).
$ws1.Select()
